$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("D3").Value = 44313
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("N3").Value = '$/caja 15 kilos empedrada'
$ws.Range("P3").Value = 1000
$ws.Range("Q3").Value = 15

# Row 4 updates
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 30000
$ws.Range("N4").Value = '$/caja 20 kilos empedrada'
$ws.Range("P4").Value = 1500
$ws.Range("Q4").Value = 20

# Row 5 updates
$ws.Range("D5").Value = 44293
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("N5").Value = '$/caja 15 kilos empedrada'
$ws.Range("P5").Value = 1667
$ws.Range("Q5").Value = 15

# Row 8 updates
$ws.Range("D8").Value = 44280
$ws.Range("J8").Value = 30
$ws.Range("N8").Value = '$/caja 18 kilos empedrada'
$ws.Range("P8").Value = 1389
$ws.Range("Q8").Value = 18
